# Apply the edit described by the diff:
#  - footer: remove the "FPOLY HCM – LƯU HÀNH NỘI BỘ" run entirely
#  - header: remove the logo picture (InlineShape) and the
#            "WEB207 – " / "Front-End Frameworks" text runs,
#            leaving just the tab run in place.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footer: drop the "FPOLY HCM – LƯU HÀNH NỘI BỘ" text ---------------
$ftr = $sec.Footers.Item(1)
$ftr.Range.Find.Execute("FPOLY HCM – LƯU HÀNH NỘI BỘ", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- Header: drop the logo image and the "WEB207 – Front-End Frameworks"
#     title text, keeping the leading tab run intact -------------------
$hdr = $sec.Headers.Item(1)

while ($hdr.Range.InlineShapes.Count -gt 0) {
    $hdr.Range.InlineShapes.Item(1).Delete()
}

$hdr.Range.Find.Execute("WEB207 – Front-End Frameworks", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 2) | Out-Null
